# Update branch to version 0.5.0
# Change supersite (column H) values:
#   Rows 12-24: "Monarch K8"  -> "Louisville MS"
#   Rows 160,162,168: "Platt MS" -> "New Horizon"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$monarchRows = 12..24
foreach ($r in $monarchRows) {
    $ws.Cells.Item($r, 8).Value = "Louisville MS"
}

$plattRows = @(160, 162, 168)
foreach ($r in $plattRows) {
    $ws.Cells.Item($r, 8).Value = "New Horizon"
}
